# Generate Report for Handoff
# Updates the "Status" to "Ready for handoff" and refreshes the handoff
# timestamps across the Overview / zh-cn / de-de sheets, and narrows the
# "Status"/"Latest Handoff Datetime" style columns that previously shared
# the wider "date" column width.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("C2").Value = "Ready for handoff"

# --- Refreshed handoff timestamps ---
$ws_overview.Range("G2").Value = "2016-08-13 17:21:15"
$ws_dede.Range("H2").Value = "2016-08-13 17:21:15"
$ws_zhcn.Range("H2").Value = "2016-08-13 17:21:07"

# --- Column width changes (Status columns narrowed to fit "Ready for handoff") ---
# Note: Excel's ColumnWidth is stored internally at whole-pixel granularity, so
# the assigned value below is the input that resolves (via Excel's
# char-width<->pixel rounding) to the narrower target column width.
$ws_overview.Range("E1").ColumnWidth = 16.3
$ws_overview.Range("F1").ColumnWidth = 16.3
$ws_zhcn.Range("C1").ColumnWidth = 16.3
$ws_dede.Range("C1").ColumnWidth = 16.3
